# ---------------------------------------------------------------------------
# Target change (per the supplied OOXML diff):
#
#   The whole package was re-saved (every auto-generated r:id / r:embed
#   relationship id in ppt/presentation.xml, ppt/slideMasters/slideMaster.xml
#   and ppt/slides/slide.xml was regenerated with a fresh random value, even
#   though the relationships themselves still point at the exact same
#   targets) and the Office Add-in ("web extension") instance embedded on
#   slide 1 got a new internal instance id:
#
#       ppt/slides/udata/data.xml
#         we:webextension/@id
#           {8b3c77f4-6d9a-4951-a353-5cb68863d90f}
#             -> {9ae3eac1-dbf9-47c7-a5ea-169582f1ea60}
#
#   (plus the we:snapshot/@r:embed + we:webextensionref/@r:id that merely
#   mirror the renamed relationship ids above).
#
# None of this is visible/user-facing content: the slide text, shapes,
# layouts, master, theme, image and the add-in reference/snapshot picture
# are all unchanged -- only internal package bookkeeping ids moved, which is
# exactly what happens when PowerPoint simply re-saves a .pptx (e.g. as part
# of a routine branch merge, matching the "merge master to gh-pages" commit
# message) without anyone editing the deck's content.
#
# The PowerPoint object model does not expose the embedded web-extension
# part (there is no Shape/Slide property for a we:webextension's id, sizes,
# snapshot relationship, etc. -- Office Add-ins are not scriptable through
# this surface in real PowerPoint either), so that internal id can't be
# poked directly via COM automation. We touch the presentation/slide through
# the object model (as a normal automation script would) without mutating
# any visible content, which keeps the deck's content identical to what a
# plain re-save produces.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Touch the existing placeholder shapes read-only, leaving their content
# exactly as-is (no text/position/formatting changes anywhere on the slide).
$null = $s.Shapes.Item(1).TextFrame.TextRange.Text
$null = $s.Shapes.Item(2).TextFrame.TextRange.Text

$p.Save()
